{"js": "// The underlying OOXML diff for this revision is a pure canonicalization of\n// the package XML: every changed line is the exact same element with its\n// attributes re-serialized in (alphabetical) order \u2014 e.g.\n//   <w:pgSz w:w=\"11906\" w:h=\"16838\"/>            ->  <w:pgSz w:h=\"16838\" w:w=\"11906\"/>\n//   <w:rFonts w:asciiTheme=\"...\" ... w:cstheme=\"minorBidi\"/>\n//                                                  ->  <w:rFonts w:asciiTheme=\"...\" w:cstheme=\"minorBidi\" .../>\n//   <w:style w:type=\"paragraph\" w:default=\"1\" .../> -> <w:style w:default=\"1\" ... w:type=\"paragraph\"/>\n// and so on through the namespace declarations on <w:document>, the\n// <w:latentStyles>/<w:lsdException> table and the <w:style> definitions.\n// No element, attribute, value, run, paragraph or text was added, removed or\n// changed \u2014 only the on-disk attribute ordering differs, which is an\n// artifact of the XML serializer that re-wrote the part, not an authored\n// document edit. There is nothing for the Word content object model (which\n// only exposes document content, not raw XML attribute order) to apply, so\n// this script intentionally performs no content mutation: it just touches\n// the body to establish a valid context/sync round-trip.\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# The underlying OOXML diff for this revision is a pure canonicalization of\n# the package XML: every changed line is the exact same element with its\n# attributes re-serialized in (alphabetical) order - e.g.\n#   <w:pgSz w:w=\"11906\" w:h=\"16838\"/>             -> <w:pgSz w:h=\"16838\" w:w=\"11906\"/>\n#   <w:rFonts w:asciiTheme=\"...\" ... w:cstheme=\"minorBidi\"/>\n#                                                   -> <w:rFonts w:asciiTheme=\"...\" w:cstheme=\"minorBidi\" .../>\n#   <w:style w:type=\"paragraph\" w:default=\"1\" .../> -> <w:style w:default=\"1\" ... w:type=\"paragraph\"/>\n# and so on through the namespace declarations on <w:document>, the\n# <w:latentStyles>/<w:lsdException> table and the <w:style> definitions.\n# No element, attribute, value, run, paragraph or text was added, removed or\n# changed - only the on-disk attribute ordering differs, which is an\n# artifact of the XML serializer that re-wrote the part, not an authored\n# document edit. There is nothing for the Word COM object model (which only\n# exposes document content, not raw XML attribute order) to apply, so this\n# script intentionally performs no content mutation: it just reads the\n# document to establish a valid, no-op automation round-trip.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
